$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '42.549.26'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '2.512.92'
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('E4').Value = '  +0.07%  '
Set-TextValue 'D5' '313.90'
$ws.Range('E5').Value = '  +4.14%  '
Set-TextValue 'D6' '95.57'
$ws.Range('E6').Value = '  -1.86%  '
Set-TextValue 'D7' '0.586'
$ws.Range('E7').Value = '  +2.35%  '
$ws.Range('E8').Value = '  +0.05%  '
Set-TextValue 'D9' '0.537'
$ws.Range('E9').Value = '  -1.02%  '
Set-TextValue 'D10' '36.44'
$ws.Range('E10').Value = '  -0.07%  '
Set-TextValue 'D11' '0.0813'
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('E13').Value = '  -2.58%  '
$ws.Range('D14').Value = '2.902.88'
$ws.Range('E14').Value = '  -1.36%  '
Set-TextValue 'D15' '15.48'
$ws.Range('E15').Value = '  +6.21%  '
$ws.Range('D16').Value = '2.525.72'
$ws.Range('E16').Value = '  -0.28%  '
Set-TextValue 'D17' '0.861'
$ws.Range('E17').Value = '  -1.76%  '
$ws.Range('D18').Value = '42.561.93'
$ws.Range('E18').Value = '  -0.26%  '
Set-TextValue 'D19' '12.83'
$ws.Range('E19').Value = '  -3.00%  '
$ws.Range('D20').Value = '0.0₃0971'
$ws.Range('E20').Value = '  -0.84%  '
Set-TextValue 'D21' '6.50'
$ws.Range('E21').Value = '  -0.73%  '
Set-TextValue 'D22' '71.52'
$ws.Range('E22').Value = '  +0.18%  '
Set-TextValue 'D23' '252.57'
$ws.Range('E23').Value = '  -0.47%  '
Set-TextValue 'D24' '2.97'
$ws.Range('E24').Value = '  +1.30%  '
Set-TextValue 'D25' '2.03'
$ws.Range('E25').Value = '  -1.38%  '
Set-TextValue 'D26' '26.98'
$ws.Range('E26').Value = '  -2.78%  '
$ws.Range('E27').Value = '  +0.04%  '
Set-TextValue 'D28' '2.34'
$ws.Range('E28').Value = '  +11.67%  '
Set-TextValue 'D29' '10.12'
$ws.Range('E29').Value = '  +1.39%  '
Set-TextValue 'D30' '37.75'
$ws.Range('E30').Value = '  -0.32%  '
Set-TextValue 'D31' '5.90'
$ws.Range('E31').Value = '  -0.87%  '
Set-TextValue 'D32' '155.15'
$ws.Range('E32').Value = '  -0.18%  '
Set-TextValue 'D33' '19.53'
$ws.Range('E33').Value = '  +6.63%  '
Set-TextValue 'D34' '3.31'
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D35' '0.0785'
$ws.Range('E35').Value = '  -1.66%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D36' '2.07'
$ws.Range('E36').Value = '  -4.64%  '
$ws.Range('E37').Value = '  -4.34%  '
$ws.Range('E38').Value = '  -0.84%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D39' '0.120'
$ws.Range('E39').Value = '  +1.44%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D40' '24.03'
$ws.Range('E40').Value = '  -6.48%  '
$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D41' '3.38'
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D42' '3.85'
$ws.Range('E42').Value = '  +0.31%  '
Set-TextValue 'D43' '2.03'
$ws.Range('E43').Value = '  -2.86%  '
Set-TextValue 'D44' '0.0302'
$ws.Range('E44').Value = '  +0.09%  '
Set-TextValue 'D45' '0.999'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '2.018.71'
$ws.Range('E46').Value = '  -2.57%  '
Set-TextValue 'D47' '84.14'
$ws.Range('E47').Value = '  -4.56%  '
Set-TextValue 'D48' '8.94'
$ws.Range('E48').Value = '  -2.80%  '
$ws.Range('D49').Value = '2.757.96'
$ws.Range('E49').Value = '  -1.48%  '
Set-TextValue 'D50' '72.99'
$ws.Range('E50').Value = '  -2.04%  '
Set-TextValue 'D51' '0.190'
$ws.Range('E51').Value = '  +0.85%  '
